$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K76").Value = 11
$ws.Range("L76").Value = 12.65

$ws.Range("K78").Value = 16.38
$ws.Range("L78").Value = 23.92

$ws.Range("K80").Value = 21.645
$ws.Range("L80").Value = 15.53

$ws.Range("K81").Value = 34.681
$ws.Range("L81").Value = 41.42

$ws.Range("K82").Value = 49.28
$ws.Range("L82").Value = 36.86

$ws.Range("K83").Value = 237.986
$ws.Range("L83").Value = 401.39
